$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 4
    5  = 3
    6  = 2
    7  = 4
    8  = 8
    9  = 5
    10 = 3
    11 = 2
    12 = 2
    13 = 4
    14 = 2
    15 = 7
    16 = 3
    17 = 4
    18 = 2
    19 = 3
    20 = 0
    21 = 5
    22 = 1
    23 = 3
    24 = 3
    25 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
